# Apply the "automatic update" changes to the Översikt SKURUP worksheet.
#
# Summary of the edit (per the diff):
#   1. Column C (row 2-9, "Förändrad" date) is bumped by one day:
#      46064 -> 46065
#   2. Rows 3 and 4 swap their A/B/G values (Beteckning, Datum, Area).
#   3. Rows 5-9 rotate their A/B/G values: each row takes on the values
#      that used to belong to the next row down, and row 9's old values
#      wrap around into row 5.
#
# Note: reading via .Value in this runtime returns a reflection
# description string instead of the actual cell value, so .Value2 is
# used for all reads/writes of plain values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" column (C) by one day for rows 2 through 9.
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}

# 2) Capture current A (Beteckning), B (Datum), G (Area) for rows 3-9
#    before overwriting anything.
$colA = @{}
$colB = @{}
$colG = @{}
for ($r = 3; $r -le 9; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

# New row -> source row mapping derived from the diff:
#   row 3 <- old row 4
#   row 4 <- old row 3
#   row 5 <- old row 6
#   row 6 <- old row 7
#   row 7 <- old row 8
#   row 8 <- old row 9
#   row 9 <- old row 5
$mapping = @{ 3 = 4; 4 = 3; 5 = 6; 6 = 7; 7 = 8; 8 = 9; 9 = 5 }

foreach ($r in @(3, 4, 5, 6, 7, 8, 9)) {
    $src = $mapping[$r]
    $ws.Cells.Item($r, 1).Value2 = $colA[$src]
    $ws.Cells.Item($r, 2).Value2 = $colB[$src]
    $ws.Cells.Item($r, 7).Value2 = $colG[$src]
}
